# Fruta / hortaliza, semanal
# The weekly refresh re-ordered the daily records: what used to be
# row 2 now belongs in row 4, and what used to be row 3 now belongs
# in row 5 (and vice versa). Swap the data-bearing columns between
# those row pairs (D, H, J, K, L, M, N, P) while leaving the columns
# that are identical across every row (A, B, C, E, F, G, I, O, Q, R)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $cols = @("D", "H", "J", "K", "L", "M", "N", "P")
    foreach ($col in $cols) {
        $cellA = $ws.Range("$col$rowA")
        $cellB = $ws.Range("$col$rowB")
        $tmp = $cellA.Value2
        $cellA.Value2 = $cellB.Value2
        $cellB.Value2 = $tmp
    }
}

Swap-Rows 2 4
Swap-Rows 3 5
